# Refresh the crypto price/volume snapshot (columns D "Price" and E
# "Volume(1h)") for rows 2-51. Rows 41/42 additionally swap which coin
# (Algorand / FraxShare) occupies which row.
#
# Price values that look like plain numbers are written with a leading
# apostrophe so Excel stores them as text, matching the workbook's existing
# text-typed Price cells instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.979.35"
$ws.Range("E2").Value = "  +5.95%  "

$ws.Range("D3").Value = "1.713.81"
$ws.Range("E3").Value = "  +3.75%  "

$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").Value = "'330.91"
$ws.Range("E5").Value = "  +5.04%  "

$ws.Range("D6").Value = "'0.9946"
$ws.Range("E6").Value = "  -0.49%  "

$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("D8").Value = "'49.02"
$ws.Range("E8").Value = "  +5.65%  "

$ws.Range("D9").Value = "'0.3331"
$ws.Range("E9").Value = "  +2.48%  "

$ws.Range("D10").Value = "'1.179"
$ws.Range("E10").Value = "  +4.85%  "

$ws.Range("D11").Value = "'0.07490"
$ws.Range("E11").Value = "  +6.70%  "

$ws.Range("D12").Value = "'0.9937"
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("E13").Value = "  +4.61%  "

$ws.Range("D14").Value = "'20.12"
$ws.Range("E14").Value = "  +3.84%  "

$ws.Range("D15").Value = "'6.897"
$ws.Range("E15").Value = "  +4.47%  "

$ws.Range("D16").Value = "1.707.76"
$ws.Range("E16").Value = "  +3.35%  "

$ws.Range("D17").Value = "'0.00001074"
$ws.Range("E17").Value = "  +3.11%  "

$ws.Range("D18").Value = "'0.06630"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("D19").Value = "'81.84"
$ws.Range("E19").Value = "  +4.10%  "

$ws.Range("D20").Value = "'0.9950"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").Value = "'16.29"
$ws.Range("E21").Value = "  +3.82%  "

$ws.Range("D22").Value = "'6.075"
$ws.Range("E22").Value = "  +2.46%  "

$ws.Range("E23").Value = "  +3.96%  "

$ws.Range("D24").Value = "25.950.01"
$ws.Range("E24").Value = "  +5.93%  "

$ws.Range("D25").Value = "'2.459"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").Value = "'2.495"
$ws.Range("E26").Value = "  +7.01%  "

$ws.Range("D27").Value = "'149.91"
$ws.Range("E27").Value = "  +1.97%  "

$ws.Range("D28").Value = "'19.27"
$ws.Range("E28").Value = "  +3.86%  "

$ws.Range("D29").Value = "'1.299"
$ws.Range("E29").Value = "  +9.27%  "

$ws.Range("D30").Value = "1.894.02"
$ws.Range("E30").Value = "  +3.28%  "

$ws.Range("D31").Value = "'129.00"
$ws.Range("E31").Value = "  +3.84%  "

$ws.Range("D32").Value = "'4.086"
$ws.Range("E32").Value = "  +0.40%  "

$ws.Range("D33").Value = "'5.971"
$ws.Range("E33").Value = "  +4.61%  "

$ws.Range("D34").Value = "'0.08510"
$ws.Range("E34").Value = "  +0.90%  "

$ws.Range("D35").Value = "'1.714"
$ws.Range("E35").Value = "  +3.17%  "

$ws.Range("E36").Value = "  +5.61%  "

$ws.Range("D37").Value = "'5.367"
$ws.Range("E37").Value = "  +3.36%  "

$ws.Range("D38").Value = "'1.286"
$ws.Range("E38").Value = "  +1.34%  "

$ws.Range("D39").Value = "'0.06213"
$ws.Range("E39").Value = "  +3.35%  "

$ws.Range("D40").Value = "'0.02286"
$ws.Range("E40").Value = "  +3.01%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.529"
$ws.Range("E41").Value = "  +5.56%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.2126"
$ws.Range("E42").Value = "  +2.81%  "

$ws.Range("D43").Value = "'14.57"
$ws.Range("E43").Value = "  +15.00%  "

$ws.Range("D44").Value = "'0.6157"
$ws.Range("E44").Value = "  +4.57%  "

$ws.Range("D45").Value = "'0.9952"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").Value = "'3.841"
$ws.Range("E46").Value = "  +1.03%  "

$ws.Range("D47").Value = "'0.5874"
$ws.Range("E47").Value = "  +4.65%  "

$ws.Range("D48").Value = "'126.85"
$ws.Range("E48").Value = "  +2.57%  "

$ws.Range("D49").Value = "'2.010"
$ws.Range("E49").Value = "  +3.46%  "

$ws.Range("D50").Value = "'0.07253"
$ws.Range("E50").Value = "  +4.80%  "

$ws.Range("D51").Value = "'76.97"
$ws.Range("E51").Value = "  +3.73%  "
